$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1237434417091066
$ws.Range("D2").Value = 0.1196842955930647
$ws.Range("E2").Value = 50.70008245461636
$ws.Range("F2").Value = 0.3108507183785727

$ws.Range("B3").Value = 0.1651129746684594
$ws.Range("D3").Value = 0.07790485870206466
$ws.Range("E3").Value = 35.71534349516163
$ws.Range("F3").Value = 0.2297039605243099

$ws.Range("B4").Value = 0.1945971317399702
$ws.Range("D4").Value = 0.1381216374723471
$ws.Range("E4").Value = 27.50637666811603
$ws.Range("F4").Value = 0.2546404718858186

$ws.Range("B5").Value = 0.2984227044273154
$ws.Range("D5").Value = 0.09791623549621251
$ws.Range("E5").Value = 23.79511536103578
$ws.Range("F5").Value = 0.2735055151383461

$ws.Range("B6").Value = 0.2575306429852456
$ws.Range("D6").Value = 0.09123370543321095
$ws.Range("E6").Value = 19.55354180051168
$ws.Range("F6").Value = 0.3310777790444115

$ws.Range("B7").Value = 0.1205744428878326
$ws.Range("D7").Value = 0.1338928226020349
$ws.Range("E7").Value = 17.89175653257555
$ws.Range("F7").Value = 0.3237905540932602

$ws.Range("B8").Value = 0.2903065877030949
$ws.Range("D8").Value = 0.1474343261998244
$ws.Range("E8").Value = 17.24357550765467
$ws.Range("F8").Value = 0.3060111889882698

$ws.Range("B9").Value = 0.2650534613228793
$ws.Range("D9").Value = 0.113436633132458
$ws.Range("E9").Value = 16.19679484058293
$ws.Range("F9").Value = 0.3176174714332895

$ws.Range("B10").Value = 0.1069799648846581
$ws.Range("D10").Value = 0.09747201301159131
$ws.Range("E10").Value = 14.73732521861051
$ws.Range("F10").Value = 0.285285616062754

$ws.Range("B11").Value = 0.1117397158165275
$ws.Range("D11").Value = 0.09511818230615347
$ws.Range("E11").Value = 14.63378371788991
$ws.Range("F11").Value = 0.2505989501941326

$ws.Range("B12").Value = 0.2389153276298857
$ws.Range("D12").Value = 0.1493404564963725
$ws.Range("E12").Value = 15.20248679319954
$ws.Range("F12").Value = 0.2920134270395738

$ws.Range("B13").Value = 0.1234063153000358
$ws.Range("D13").Value = 0.08172077756908988
$ws.Range("E13").Value = 13.48536236530232
$ws.Range("F13").Value = 0.2210044451840307

$ws.Range("B14").Value = 0.1637043449097016
$ws.Range("D14").Value = 0.1211847624708794
$ws.Range("E14").Value = 13.26973472070217
$ws.Range("F14").Value = 0.3739776655641557

$ws.Range("B15").Value = 0.1416137195523914
$ws.Range("D15").Value = 0.05858195690155565
$ws.Range("E15").Value = 14.18484577557363
$ws.Range("F15").Value = 0.3382202091551558

$ws.Range("B16").Value = 0.1657415923017283
$ws.Range("D16").Value = 0.1236512837913234
$ws.Range("E16").Value = 14.01144982241755
$ws.Range("F16").Value = 0.3954414105569229

$ws.Range("B17").Value = 0.1521030928814036
$ws.Range("D17").Value = 0.06395214164237525
$ws.Range("E17").Value = 14.26936307540018
$ws.Range("F17").Value = 0.3699104877016595

$ws.Range("B18").Value = 0.2156421199695095
$ws.Range("D18").Value = 0.09664268266212483
$ws.Range("E18").Value = 13.40450800653233
$ws.Range("F18").Value = 0.2314518626105827

$ws.Range("B19").Value = 0.1520633567910731
$ws.Range("D19").Value = 0.09365213748078141
$ws.Range("E19").Value = 14.42913596613205
$ws.Range("F19").Value = 0.3951883294469988

$ws.Range("B20").Value = 0.2250188688028919
$ws.Range("D20").Value = 0.1481841363105223
$ws.Range("E20").Value = 14.6713028939412
$ws.Range("F20").Value = 0.2135786459926527

$ws.Range("B21").Value = 0.1128081393942871
$ws.Range("D21").Value = 0.09446184442918287
$ws.Range("E21").Value = 14.99892619187818
$ws.Range("F21").Value = 0.2407511775494566
